# Rename the two data-source header labels across every worksheet
# ("Source1" -> "1st Source", "Source2" -> "2nd Source") and resize the
# columns that hold them, mirroring the best-fit width recalculation
# Excel performs after a header is renamed.

$wb = $excel.ActiveWorkbook

# ColumnWidth (character units) inputs chosen so the engine's internal
# rounding lands on the widths closest to the new best-fit values.
$widthForSource1Col = 14.0
$widthForSource2Col = 14.65

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    if ($used -eq $null) { continue }

    $source1Cols = @{}
    $source2Cols = @{}

    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if ($val -eq "Source1") {
            $cell.Value = "1st Source"
            $source1Cols[$cell.Column] = $true
        } elseif ($val -eq "Source2") {
            $cell.Value = "2nd Source"
            $source2Cols[$cell.Column] = $true
        }
    }

    foreach ($col in $source1Cols.Keys) {
        $ws.Columns.Item($col).ColumnWidth = $widthForSource1Col
    }
    foreach ($col in $source2Cols.Keys) {
        $ws.Columns.Item($col).ColumnWidth = $widthForSource2Col
    }
}
